# V 2.0.2 se arreglo la fechar y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, avoiding Excel's automatic
# number/date reinterpretation (used for purely-numeric / date-like strings)
# by routing the write through a text formula -> Copy -> PasteSpecial(values)
# round-trip, which preserves the cell's existing style.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Patient name / expediente clinico
$ws.Range("A6").Value = "COSAJAY  MEJIA  JOSSELYN   ESTEFANY"
$ws.Range("G6").Value = "8257/201761502"

# Fecha de nacimiento / Edad / Lugar de nacimiento
Set-TextValue $ws.Range("A9") "1994-06-04"
$ws.Range("D9").Value = "23 AÑOS "
$ws.Range("E9").Value = "GUATEMALA"

# Estado civil
$ws.Range("A11").Value = "SOLTERO"

# Documento de identificacion
Set-TextValue $ws.Range("G11") "2434796970101"

# Contacto de emergencia
$ws.Range("A13").Value = "PATRICIA MEJIA"
$ws.Range("D13").Value = "MAMA"
$ws.Range("E13").Value = "35 AV 10-41 Z18 PARAISO 2"
Set-TextValue $ws.Range("G13") "55007613"

# Fecha / hora de la asistencia medica
$ws.Range("D14").Value = "Hora: 17:34:41"
$ws.Range("A15").Value = "22/10/2017"

$excel.CutCopyMode = $false
